# Matrix Algebra Slides Extra - "Make new figure for Markov tunnels"
#
# 1. Slide 10 (TextBox 11 / "In R" box): move/resize to sit above the new
#    (narrower) picture.
# 2. Slide 2 ("Estimation of proportion..." textbox): reword first
#    paragraph, and trim the last formula paragraph so it stops right
#    after "Pr(Dead)".
# 3. Slide 9: highlight the "*" and "%*%" operators in accent5 colour.
# 4. Delete the trailing empty slide (slide 11).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 10 - resize/reposition the "In R" label box next to the picture
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
for ($i = 1; $i -le $s10.Shapes.Count; $i++) {
    $shp = $s10.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 11") {
        $shp.Left = 396.46181102362203
        $shp.Width = 260.5382677165354
    }
}

# ---------------------------------------------------------------------
# 2. Slide 2 - edit the proportion-estimation textbox
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame -eq $false) { continue }
    $tr = $shp.TextFrame.TextRange
    if (-not $tr.Text.Contains("Estimation of proportion")) { continue }

    # 2a. First paragraph: "... proportion in a specific ..." ->
    #     "... proportion of a cohort in a specific ..."
    $paraCount = $tr.Paragraphs().Count
    for ($i2 = 1; $i2 -le $paraCount; $i2++) {
        $para = $tr.Paragraphs($i2)
        if ($para.Text.StartsWith("Estimation of proportion")) {
            $c = $tr.Characters($para.Start, $para.Length)
            $c.Text = "Estimation of proportion of a cohort in a specific health states also involves multiplication and addition."
        }
    }

    # 2b. Last formula paragraph: cut everything from " - (psick, dead )*Pr(Sick)"
    #     so the line ends right after the (italic) "Dead" run, keeping the
    #     closing ")" run intact.
    $paraCount = $tr.Paragraphs().Count
    for ($i2 = 1; $i2 -le $paraCount; $i2++) {
        $para = $tr.Paragraphs($i2)
        if ($para.Text.Contains("pDead") -and $para.Text.Contains("Pr(Sick)")) {
            $n = $tr.Runs().Count
            $deadRunEndPos = -1
            $closeParenStart = -1
            for ($j = 1; $j -le $n; $j++) {
                $r = $tr.Runs($j)
                if ($r.Start -ge $para.Start -and $r.Start -lt ($para.Start + $para.Length)) {
                    if ($r.Text -eq "Dead") {
                        $deadRunEndPos = $r.Start + 4
                    }
                    if ($r.Text -eq ")" -and $deadRunEndPos -ne -1) {
                        $closeParenStart = $r.Start
                    }
                }
            }
            if ($deadRunEndPos -ne -1 -and $closeParenStart -ne -1 -and $closeParenStart -gt $deadRunEndPos) {
                $delLen = $closeParenStart - $deadRunEndPos
                $delRange = $tr.Characters($deadRunEndPos, $delLen)
                $delRange.Text = ""
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3. Slide 9 - colour the "*" and "%*%" operators with accent5
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
for ($i = 1; $i -le $s9.Shapes.Count; $i++) {
    $shp = $s9.Shapes.Item($i)
    if ($shp.HasTextFrame -eq $false) { continue }
    $tr = $shp.TextFrame.TextRange
    $fullText = $tr.Text

    if ($fullText.Contains("multiplication operator in R") -and $fullText.Contains("gives")) {
        $r1 = $tr.Runs(1)
        $rel = $r1.Text.IndexOf("*")
        if ($rel -ge 0) {
            $target = $tr.Characters($r1.Start + $rel, 1)
            $target.Font.Color.ObjectThemeColor = 9
        }
    }

    if ($fullText.Contains("achieved using the") -and $fullText.Contains("%*%")) {
        $r1 = $tr.Runs(1)
        $rel = $r1.Text.IndexOf("%*%")
        if ($rel -ge 0) {
            $target = $tr.Characters($r1.Start + $rel, 3)
            $target.Font.Color.ObjectThemeColor = 9
        }
    }
}

# ---------------------------------------------------------------------
# 4. Remove the trailing empty slide (slide 11)
# ---------------------------------------------------------------------
$p.Slides.Item($p.Slides.Count).Delete()
